# Daily attendance processing - 2026-01-16 10:36:21
#
# The "Recorded By" column (G) lists session recorders as a
# comma-separated string. Swap the order of the two names so
# "System, dnasr281@gmail.com" becomes "dnasr281@gmail.com, System"
# for every row in the used range that currently holds that value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = $cell.Value2
    if ($current -eq $oldValue) {
        $cell.Value = $newValue
    }
}
